{"js": "// The diff replaces the text of each arithmetic-expression cell in the\n// document's single table, in row-major (reading) order: row 0 cell 0,\n// row 0 cell 1, ... row 0 cell 4, row 1 cell 0, ... This array holds the\n// old/new text pairs in that exact order (old text is only used for a\n// sanity check; the replacement always happens by position).\nconst REPLACEMENTS = [\n  [\"69+8=77\", \"98-59=39\"],\n  [\"11+60=71\", \"81-38=43\"],\n  [\"14+82=96\", \"41+20=61\"],\n  [\"9+7=16\", \"59-23=36\"],\n  [\"5+34=39\", \"39+34=73\"],\n  [\"69+15=84\", \"95-67=28\"],\n  [\"60-56=4\", \"75-48=27\"],\n  [\"26-13=13\", \"77-58=19\"],\n  [\"74-72=2\", \"55-1=54\"],\n  [\"32-27=5\", \"98-4=94\"],\n  [\"42+48=90\", \"89-31=58\"],\n  [\"49-48=1\", \"92-88=4\"],\n  [\"64-38=26\", \"44-13=31\"],\n  [\"31-20=11\", \"81+7=88\"],\n  [\"33+35=68\", \"32+0=32\"],\n  [\"54+40=94\", \"80-46=34\"],\n  [\"41-11=30\", \"96-11=85\"],\n  [\"18+3=21\", \"12+60=72\"],\n  [\"33-25=8\", \"37+43=80\"],\n  [\"57+10=67\", \"23-20=3\"],\n  [\"69-41=28\", \"15+0=15\"],\n  [\"82-66=16\", \"61-24=37\"],\n  [\"75-4=71\", \"40+1=41\"],\n  [\"17+4=21\", \"73+7=80\"],\n  [\"28+17=45\", \"68-54=14\"],\n  [\"55+12=67\", \"67-42=25\"],\n  [\"18+4=22\", \"1+24=25\"],\n  [\"2+65=67\", \"62-3=59\"],\n  [\"28-21=7\", \"1+83=84\"],\n  [\"95-57=38\", \"69+27=96\"],\n  [\"62+9=71\", \"11+76=87\"],\n  [\"18+43=61\", \"63-35=28\"],\n  [\"70-61=9\", \"49+10=59\"],\n  [\"44+12=56\", \"99-18=81\"],\n  [\"8+86=94\", \"42+39=81\"],\n  [\"86-23=63\", \"34+58=92\"],\n  [\"77-42=35\", \"54-34=20\"],\n  [\"9+69=78\", \"9+10=19\"],\n  [\"64-36=28\", \"89+10=99\"],\n  [\"16+53=69\", \"3+54=57\"],\n  [\"74+20=94\", \"97-4=93\"],\n  [\"97-74=23\", \"29+66=95\"],\n  [\"84-40=44\", \"98-98=0\"],\n  [\"18+71=89\", \"4+49=53\"],\n  [\"86-23=63\", \"73-10=63\"],\n  [\"19+29=48\", \"9+87=96\"],\n  [\"38+36=74\", \"94-73=21\"],\n  [\"95-39=56\", \"51+10=61\"],\n  [\"3+47=50\", \"2+71=73\"],\n  [\"66-27=39\", \"98-30=68\"],\n  [\"37+46=83\", \"66-12=54\"],\n  [\"41-28=13\", \"61+24=85\"],\n  [\"0+1=1\", \"59-9=50\"],\n  [\"57+13=70\", \"8+90=98\"],\n  [\"89-55=34\", \"32-6=26\"],\n  [\"84+14=98\", \"27-5=22\"],\n  [\"87-16=71\", \"60-43=17\"],\n  [\"12+36=48\", \"18+39=57\"],\n  [\"44-20=24\", \"62+1=63\"],\n  [\"45+10=55\", \"37-33=4\"],\n  [\"55+10=65\", \"26+35=61\"],\n  [\"23-5=18\", \"25+21=46\"],\n  [\"92-44=48\", \"72-8=64\"],\n  [\"30+61=91\", \"15+39=54\"],\n  [\"38-34=4\", \"41-22=19\"],\n  [\"37+13=50\", \"59+36=95\"],\n  [\"80-31=49\", \"17+74=91\"],\n  [\"44+6=50\", \"82-41=41\"],\n  [\"68-48=20\", \"33+23=56\"],\n  [\"8+68=76\", \"92-65=27\"],\n  [\"67-47=20\", \"0+20=20\"],\n  [\"67-26=41\", \"90-63=27\"],\n  [\"20+10=30\", \"95-44=51\"],\n  [\"4+18=22\", \"13+42=55\"],\n  [\"79+17=96\", \"25+65=90\"],\n  [\"4+61=65\", \"44-32=12\"],\n  [\"76-65=11\", \"22+12=34\"],\n  [\"42+0=42\", \"7+12=19\"],\n  [\"19-7=12\", \"84-34=50\"],\n  [\"92-39=53\", \"57+33=90\"],\n  [\"24-23=1\", \"24+66=90\"],\n  [\"50+43=93\", \"10+4=14\"],\n  [\"16+45=61\", \"80-25=55\"],\n  [\"36-13=23\", \"94-32=62\"],\n  [\"12+23=35\", \"78-29=49\"],\n  [\"78-13=65\", \"97-54=43\"],\n  [\"10+17=27\", \"57+41=98\"],\n  [\"27-23=4\", \"71-37=34\"],\n  [\"89-38=51\", \"30+17=47\"],\n  [\"34+3=37\", \"13+61=74\"],\n  [\"87-73=14\", \"61-25=36\"],\n  [\"84-8=76\", \"10+89=99\"],\n  [\"95-11=84\", \"15+21=36\"],\n  [\"30+28=58\", \"41+17=58\"],\n  [\"80-40=40\", \"4+89=93\"],\n  [\"18+76=94\", \"1+57=58\"],\n  [\"39-1=38\", \"4+53=57\"],\n  [\"51+16=67\", \"51+15=66\"],\n  [\"58-2=56\", \"71-70=1\"],\n  [\"83-39=44\", \"46+0=46\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nfor (const row of rows) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Flatten all cells in row-major order, load current values.\nconst allCells = [];\nfor (const row of rows) {\n  for (const cell of row.cells.items) {\n    allCells.push(cell);\n    cell.load(\"value\");\n  }\n}\nawait context.sync();\n\nif (allCells.length !== REPLACEMENTS.length) {\n  throw new Error(\n    \"Cell count (\" + allCells.length + \") does not match replacement count (\" +\n      REPLACEMENTS.length + \").\"\n  );\n}\n\nfor (let i = 0; i < allCells.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const cell = allCells[i];\n  const current = (cell.value || \"\").trim();\n  if (current !== oldText) {\n    throw new Error(\n      \"Cell \" + i + \": expected \\\"\" + oldText + \"\\\" but found \\\"\" + current + \"\\\".\"\n    );\n  }\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# The diff replaces the text of each arithmetic-expression cell in the\n# document's single table, in row-major (reading) order: row 1 cell 1,\n# row 1 cell 2, ... row 1 cell 5, row 2 cell 1, ... This array holds the\n# old/new text pairs in that exact order (old text is only used for a\n# sanity check; the replacement always happens by position).\n$replacements = @(\n    @(\"69+8=77\", \"98-59=39\"),\n    @(\"11+60=71\", \"81-38=43\"),\n    @(\"14+82=96\", \"41+20=61\"),\n    @(\"9+7=16\", \"59-23=36\"),\n    @(\"5+34=39\", \"39+34=73\"),\n    @(\"69+15=84\", \"95-67=28\"),\n    @(\"60-56=4\", \"75-48=27\"),\n    @(\"26-13=13\", \"77-58=19\"),\n    @(\"74-72=2\", \"55-1=54\"),\n    @(\"32-27=5\", \"98-4=94\"),\n    @(\"42+48=90\", \"89-31=58\"),\n    @(\"49-48=1\", \"92-88=4\"),\n    @(\"64-38=26\", \"44-13=31\"),\n    @(\"31-20=11\", \"81+7=88\"),\n    @(\"33+35=68\", \"32+0=32\"),\n    @(\"54+40=94\", \"80-46=34\"),\n    @(\"41-11=30\", \"96-11=85\"),\n    @(\"18+3=21\", \"12+60=72\"),\n    @(\"33-25=8\", \"37+43=80\"),\n    @(\"57+10=67\", \"23-20=3\"),\n    @(\"69-41=28\", \"15+0=15\"),\n    @(\"82-66=16\", \"61-24=37\"),\n    @(\"75-4=71\", \"40+1=41\"),\n    @(\"17+4=21\", \"73+7=80\"),\n    @(\"28+17=45\", \"68-54=14\"),\n    @(\"55+12=67\", \"67-42=25\"),\n    @(\"18+4=22\", \"1+24=25\"),\n    @(\"2+65=67\", \"62-3=59\"),\n    @(\"28-21=7\", \"1+83=84\"),\n    @(\"95-57=38\", \"69+27=96\"),\n    @(\"62+9=71\", \"11+76=87\"),\n    @(\"18+43=61\", \"63-35=28\"),\n    @(\"70-61=9\", \"49+10=59\"),\n    @(\"44+12=56\", \"99-18=81\"),\n    @(\"8+86=94\", \"42+39=81\"),\n    @(\"86-23=63\", \"34+58=92\"),\n    @(\"77-42=35\", \"54-34=20\"),\n    @(\"9+69=78\", \"9+10=19\"),\n    @(\"64-36=28\", \"89+10=99\"),\n    @(\"16+53=69\", \"3+54=57\"),\n    @(\"74+20=94\", \"97-4=93\"),\n    @(\"97-74=23\", \"29+66=95\"),\n    @(\"84-40=44\", \"98-98=0\"),\n    @(\"18+71=89\", \"4+49=53\"),\n    @(\"86-23=63\", \"73-10=63\"),\n    @(\"19+29=48\", \"9+87=96\"),\n    @(\"38+36=74\", \"94-73=21\"),\n    @(\"95-39=56\", \"51+10=61\"),\n    @(\"3+47=50\", \"2+71=73\"),\n    @(\"66-27=39\", \"98-30=68\"),\n    @(\"37+46=83\", \"66-12=54\"),\n    @(\"41-28=13\", \"61+24=85\"),\n    @(\"0+1=1\", \"59-9=50\"),\n    @(\"57+13=70\", \"8+90=98\"),\n    @(\"89-55=34\", \"32-6=26\"),\n    @(\"84+14=98\", \"27-5=22\"),\n    @(\"87-16=71\", \"60-43=17\"),\n    @(\"12+36=48\", \"18+39=57\"),\n    @(\"44-20=24\", \"62+1=63\"),\n    @(\"45+10=55\", \"37-33=4\"),\n    @(\"55+10=65\", \"26+35=61\"),\n    @(\"23-5=18\", \"25+21=46\"),\n    @(\"92-44=48\", \"72-8=64\"),\n    @(\"30+61=91\", \"15+39=54\"),\n    @(\"38-34=4\", \"41-22=19\"),\n    @(\"37+13=50\", \"59+36=95\"),\n    @(\"80-31=49\", \"17+74=91\"),\n    @(\"44+6=50\", \"82-41=41\"),\n    @(\"68-48=20\", \"33+23=56\"),\n    @(\"8+68=76\", \"92-65=27\"),\n    @(\"67-47=20\", \"0+20=20\"),\n    @(\"67-26=41\", \"90-63=27\"),\n    @(\"20+10=30\", \"95-44=51\"),\n    @(\"4+18=22\", \"13+42=55\"),\n    @(\"79+17=96\", \"25+65=90\"),\n    @(\"4+61=65\", \"44-32=12\"),\n    @(\"76-65=11\", \"22+12=34\"),\n    @(\"42+0=42\", \"7+12=19\"),\n    @(\"19-7=12\", \"84-34=50\"),\n    @(\"92-39=53\", \"57+33=90\"),\n    @(\"24-23=1\", \"24+66=90\"),\n    @(\"50+43=93\", \"10+4=14\"),\n    @(\"16+45=61\", \"80-25=55\"),\n    @(\"36-13=23\", \"94-32=62\"),\n    @(\"12+23=35\", \"78-29=49\"),\n    @(\"78-13=65\", \"97-54=43\"),\n    @(\"10+17=27\", \"57+41=98\"),\n    @(\"27-23=4\", \"71-37=34\"),\n    @(\"89-38=51\", \"30+17=47\"),\n    @(\"34+3=37\", \"13+61=74\"),\n    @(\"87-73=14\", \"61-25=36\"),\n    @(\"84-8=76\", \"10+89=99\"),\n    @(\"95-11=84\", \"15+21=36\"),\n    @(\"30+28=58\", \"41+17=58\"),\n    @(\"80-40=40\", \"4+89=93\"),\n    @(\"18+76=94\", \"1+57=58\"),\n    @(\"39-1=38\", \"4+53=57\"),\n    @(\"51+16=67\", \"51+15=66\"),\n    @(\"58-2=56\", \"71-70=1\"),\n    @(\"83-39=44\", \"46+0=46\")\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\nif (($rowCount * $colCount) -ne $replacements.Count) {\n    throw (\"Cell count (\" + ($rowCount * $colCount) + \") does not match replacement count (\" + $replacements.Count + \").\")\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $pair = $replacements[$i]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n        $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($current -ne $oldText) {\n            throw (\"Cell \" + $i + \" (row \" + $r + \", col \" + $c + \"): expected [\" + $oldText + \"] but found [\" + $current + \"].\")\n        }\n        $cell.Range.Text = $newText\n        $i++\n    }\n}\n\nWrite-Output (\"Updated \" + $i + \" cells.\")\n"}
